# Adds two new product rows into the CNA01D worksheet (Sheet1).
#
# 1) Insert "20140419 / IDM CD WNT SOFT 2-XL" right after the
#    "IDM CD WN.BMB2S L/XL" row -- this becomes the new row 11, and every
#    row at/below the old row 11 shifts down by one.
# 2) Insert "20140363 / LRST ANKLE MSTY 1+1" right after the
#    "LARIST KK ANKLE SOCK" row -- once the first insert above has shifted
#    things down, that row is now row 33, so the new row lands at row 34,
#    and every row at/below shifts down by one more.
#
# Both new rows are given the same look (thin border / General number
# format) as the surrounding data rows, and their values are written as
# text -- matching how every other cell in the sheet (even purely
# numeric-looking ones like "1", "2", "PT") is stored as a shared string
# rather than a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 11: new SKU 20140419 --------------------------------------------
$ws.Rows("11:11").Insert()

$dst1 = $ws.Range("A11:F11")
$dst1.NumberFormat = "@"
$ws.Range("A11").Value = "20140419"
$ws.Range("B11").Value = "IDM CD WNT SOFT 2-XL"
$ws.Range("C11").Value = "CNA01D"
$ws.Range("D11").Value = "2"
$ws.Range("E11").Value = "2"
$ws.Range("F11").Value = "PT"

# Re-apply the formatting (thin border, General number format, column
# widths, ...) used by the row right below, without touching the values
# just written.
$src1 = $ws.Range("A12:F12")
$src1.Copy()
$dst1.PasteSpecial(-4122)

# ---- Row 34: new SKU 20140363 --------------------------------------------
$ws.Rows("34:34").Insert()

$dst2 = $ws.Range("A34:F34")
$dst2.NumberFormat = "@"
$ws.Range("A34").Value = "20140363"
$ws.Range("B34").Value = "LRST ANKLE MSTY 1+1"
$ws.Range("C34").Value = "CNA01D"
$ws.Range("D34").Value = "4"
$ws.Range("E34").Value = "8"
$ws.Range("F34").Value = "RT"

$src2 = $ws.Range("A35:F35")
$src2.Copy()
$dst2.PasteSpecial(-4122)
